$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (MARS, SVM) ---
$ws.Range("A6").Value = "MARS"
$ws.Range("B6").Value = 1.1499889999999999
$ws.Range("C6").Value = 0.12764130000000001

$ws.Range("A7").Value = "SVM"
$ws.Range("B7").Value = 1.181127
$ws.Range("C7").Value = 0.1332045

# --- Number formatting: accounting-style format with 4 decimals applied
#     to the header row (bold) and the whole numeric data block ---
$acctFormat = "_(* #,##0.0000_);_(* \(#,##0.0000\);_(* ""-""??_);_(@_)"
$ws.Range("B1:C1").NumberFormat = $acctFormat
$ws.Range("B2:C7").NumberFormat = $acctFormat

# --- Column widths (approximate best-fit sizing for the new layout) ---
$ws.Columns("A").ColumnWidth = 9.43
$ws.Range("B:C").ColumnWidth = 8.3

# --- Selection, matching the saved cursor position in the workbook ---
$ws.Range("E12").Select() | Out-Null
